# Append two new feed rows to the "Filtered Feeds" sheet, matching the
# upstream workflow's new scrape results (Promega MSI Detection Kit / China
# NMPA approval as CDx for Merck's Keytruda).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$title = "Promega MSI Detection Kit Gets China NMPA Approval as CDx for Merck's Keytruda"

# Row 92: GenomeWeb article
$ws.Range("B92").Value() = "CDx"
$ws.Range("C92").Value() = $title
$ws.Hyperlinks.Add($ws.Range("A92"), "https://www.genomeweb.com/cancer/promega-msi-detection-kit-gets-china-nmpa-approval-cdx-mercks-keytruda", "", "", "https://www.genomeweb.com/cancer/promega-msi-detection-kit-gets-china-nmpa-approval-cdx-mercks-keytruda")
$ws.Range("A92").Style = "Hyperlink"

# Row 93: 360Dx article (same headline)
$ws.Range("B93").Value() = "CDx"
$ws.Range("C93").Value() = $title
$ws.Hyperlinks.Add($ws.Range("A93"), "https://www.360dx.com/cancer/promega-msi-detection-kit-gets-china-nmpa-approval-cdx-mercks-keytruda", "", "", "https://www.360dx.com/cancer/promega-msi-detection-kit-gets-china-nmpa-approval-cdx-mercks-keytruda")
$ws.Range("A93").Style = "Hyperlink"
